$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.174.14'
$ws.Range('E2').Value = '  +1.00%  '

$ws.Range('D3').Value = '1.572.12'
$ws.Range('E3').Value = '  +1.33%  '

$ws.Range('D4').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D4').Value = '1.02'
$ws.Range('E4').Value = '  +1.60%  '

$ws.Range('D5').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D5').Value = '211.39'
$ws.Range('E5').Value = '  +2.50%  '

$ws.Range('E6').Value = '  +0.98%  '

$ws.Range('E7').Value = '  +1.06%  '

$ws.Range('D8').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D8').Value = '22.08'
$ws.Range('E8').Value = '  +0.59%  '

$ws.Range('D9').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D9').Value = '0.249'
$ws.Range('E9').Value = '  +0.81%  '

$ws.Range('E10').Value = '  +0.81%  '

$ws.Range('D11').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D11').Value = '0.0865'
$ws.Range('E11').Value = '  +1.03%  '

$ws.Range('D12').Value = '1.791.59'
$ws.Range('E12').Value = '  +1.05%  '

$ws.Range('D13').Value = '1.573.81'
$ws.Range('E13').Value = '  +1.42%  '

$ws.Range('D14').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D14').Value = '3.78'
$ws.Range('E14').Value = '  +0.87%  '

$ws.Range('E15').Value = '  +0.38%  '

$ws.Range('D16').Value = '27.180.21'
$ws.Range('E16').Value = '  +1.00%  '

$ws.Range('D17').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D17').Value = '62.32'
$ws.Range('E17').Value = '  +1.18%  '

$ws.Range('E18').Value = '  -1.05%  '

$ws.Range('D19').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D19').Value = '215.97'
$ws.Range('E19').Value = '  -0.55%  '

$ws.Range('D20').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D20').Value = '7.39'
$ws.Range('E20').Value = '  +1.37%  '

$ws.Range('E21').Value = '  +1.21%  '

$ws.Range('D22').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D22').Value = '4.15'
$ws.Range('E22').Value = '  +1.77%  '

$ws.Range('D23').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D23').Value = '9.21'
$ws.Range('E23').Value = '  +0.37%  '

$ws.Range('E24').Value = '  +0.71%  '

$ws.Range('D25').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D25').Value = '154.55'
$ws.Range('E25').Value = '  +0.66%  '

$ws.Range('E26').Value = '  -0.39%  '

$ws.Range('D27').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D27').Value = '15.13'
$ws.Range('E27').Value = '  +0.98%  '

$ws.Range('E28').Value = '  +1.61%  '

$ws.Range('E29').Value = '  +1.24%  '

$ws.Range('E30').Value = '  +5.50%  '

$ws.Range('D31').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D31').Value = '0.0474'
$ws.Range('E31').Value = '  +1.22%  '

$ws.Range('D32').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D32').Value = '3.25'
$ws.Range('E32').Value = '  +1.14%  '

$ws.Range('E33').Value = '  +2.63%  '

$ws.Range('D34').Value = '1.434.75'
$ws.Range('E34').Value = '  +1.60%  '

$ws.Range('D35').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D35').Value = '1.10'
$ws.Range('E35').Value = '  +13.46%  '

$ws.Range('E36').Value = '  +1.17%  '

$ws.Range('D37').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D37').Value = '2.37'
$ws.Range('E37').Value = '  +3.36%  '

$ws.Range('E38').Value = '  +1.18%  '

$ws.Range('E39').Value = '  +1.56%  '

$ws.Range('D40').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.44'
$ws.Range('E40').Value = '  +5.93%  '

$ws.Range('D41').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '5.86'
$ws.Range('E41').Value = '  +3.63%  '

$ws.Range('D42').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '0.811'
$ws.Range('E42').Value = '  +0.46%  '

$ws.Range('D43').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '1.01'
$ws.Range('E43').Value = '  +1.37%  '

$ws.Range('D44').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D44').Value = '1.01'
$ws.Range('E44').Value = '  +0.98%  '

$ws.Range('D45').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D45').Value = '64.75'
$ws.Range('E45').Value = '  +0.49%  '

$ws.Range('E46').Value = '  +0.71%  '

$ws.Range('D47').Value = '1.710.80'

$ws.Range('D48').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D48').Value = '85.74'
$ws.Range('E48').Value = '  -1.63%  '

$ws.Range('D49').Value = '0.0₆0102'
$ws.Range('E49').Value = '  +1.61%  '

$ws.Range('D50').NumberFormat = '@'  # keep numeric-looking price as text
$ws.Range('D50').Value = '0.0518'
$ws.Range('E50').Value = '  -0.60%  '

$ws.Range('E51').Value = '  +0.32%  '
